$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.477.62"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.87%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.871.93"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.60%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.001"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.09%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'0.7182"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.72%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'239.25"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.56%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.15%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.07810"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -4.76%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.3071"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.84%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'25.24"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +8.77%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.08240"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.89%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'1.870.22"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +0.09%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'5.229"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +1.03%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'0.7209"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.70%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'89.99"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.35%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'29.531.22"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.08%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'5.828"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.66%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.21%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'240.84"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.48%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.52%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'2.121.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.05%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'1.001"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.04%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'1.001"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.08%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'7.729"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +4.03%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.1564"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +7.13%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'162.66"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.00%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'8.978"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.20%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'18.32"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +1.26%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D30').Value = "'1.358"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -4.70%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.14%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'4.328"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -1.74%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'4.078"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.29%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.05249"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.48%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'1.198"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.41%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.7163"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +1.11%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'1.000"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.11%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'2.674"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.10%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.01868"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.46%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'2.720"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.34%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'1.174.50"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +3.04%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.9058"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -1.97%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'5.992"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +1.58%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.4307"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.51%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'71.36"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +1.62%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.20%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.59%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'0.5358"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.66%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.80%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'9.144"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.74%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'7.024"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.98%  "
$ws.Range('E51').Style = 'Normal'
